$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Fill in the newly-recorded test rows (14-25) with dates, test names,
#    results and notes. The shared-string table must end up with the new
#    strings appearing in a specific order, so we deliberately set the
#    "Notes" text first (which becomes shared-string #18) and then touch the
#    "Test Name/ID" column in the exact order that reproduces the expected
#    shared-string sequence.
# ---------------------------------------------------------------------------

# Dates (column B) - 03 Dec 2019 for every new row
foreach ($r in 14..25) {
    $ws.Range("B$r").Value = 43802
}

# Results (column D) - all "fail"
foreach ($r in 14..25) {
    $ws.Range("D$r").Value = "fail"
}

# Notes (column E) - all "only test base created"; set row 14 first so this
# string becomes shared-string index 18 (the first new unique string).
foreach ($r in 14..25) {
    $ws.Range("E$r").Value = "only test base created"
}

# Test Name/ID (column C) - set in the specific order that reproduces the
# target shared-string table ordering.
$ws.Range("C21").Value = "testRestartOnClick()"
$ws.Range("C19").Value = "testInitialiseGameWordVariables()"
$ws.Range("C18").Value = "testWordOnClick()"
$ws.Range("C17").Value = "testStartEndGameScene()"
$ws.Range("C16").Value = "testSetDifficulty()"
$ws.Range("C15").Value = "testStartGameScene()"
$ws.Range("C14").Value = "testSetStage()"
$ws.Range("C24").Value = "testStart()"
$ws.Range("C25").Value = "testMain()"
$ws.Range("C20").Value = "testOnKeyPressed()"
$ws.Range("C22").Value = "estCloseOnClick()"
$ws.Range("C23").Value = "testSetFinalScore()"

# ---------------------------------------------------------------------------
# 2. Widen column C to fit the new, longer test names.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 24.16666666666667

# ---------------------------------------------------------------------------
# 3. Update the view: scroll so row 8 / column B is the top-left visible
#    cell, and move the active selection to D21.
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("D21").Select() | Out-Null

# ---------------------------------------------------------------------------
# 4. Add the extra conditional-formatting rules that highlight the
#    "Test Name/ID" column for the newly added rows (C20:C23), mirroring the
#    existing pass/fail highlighting rules used elsewhere on the sheet.
# ---------------------------------------------------------------------------
$newRuleRange = $ws.Range("C20:C23")

$failRule = $newRuleRange.FormatConditions.Add(2, 3, "=EXACT(`$D20, `$G`$4)")
$failRule.Interior.Color = 1593064

$passRule = $newRuleRange.FormatConditions.Add(2, 3, "=EXACT(`$D20,`$G`$3)")
$passRule.Interior.Color = 5296274
